$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.468607902526855
$ws.Range("B1").Value = 3.309092998504639
$ws.Range("C1").Value = 2.544580698013306
$ws.Range("D1").Value = 2.352552175521851
$ws.Range("E1").Value = 1.964970111846924
